# Applies the "update new orleans xlsx files" edit:
#   1. hotel_info gains a new "State" column (value "Louisiana") inserted
#      right after "Hotel_Name" and before "City".
#   2. The sheet order is swapped so "review_info" becomes the first tab
#      and "hotel_info" becomes the second tab.

$wb = $excel.ActiveWorkbook

$hotelSheet = $wb.Worksheets.Item("hotel_info")
$reviewSheet = $wb.Worksheets.Item("review_info")

# Insert a new column C (State) in hotel_info, shifting City/Zip/... right.
$hotelSheet.Columns.Item(3).Insert()
$hotelSheet.Cells.Item(1, 3).Value = "State"
$hotelSheet.Cells.Item(2, 3).Value = "Louisiana"

# Reorder the tabs: review_info first, hotel_info second.
$reviewSheet.Move($wb.Worksheets.Item(1))
